# ---------------------------------------------------------------------------
# feat: add 2022-Q4 data
#
# 1. Insert a new "2022-Q4" sheet right after "总计", pushing the existing
#    "2022-Q3" and "2022-Q1" sheets one slot to the right.
# 2. Populate "2022-Q4" with the three new fund-holding rows.
# 3. Update the "总计" summary sheet with a new row for 2022-Q4 (and keep
#    the existing 2022-Q3 / 2022-Q1 summary rows, renumbering the running
#    index column).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- locate the existing sheets (by position, before any insert) ----------
$summary = $wb.Worksheets.Item(1)   # "总计"
$q3      = $wb.Worksheets.Item(2)   # "2022-Q3"  (currently 2nd sheet)
$q1      = $wb.Worksheets.Item(3)   # "2022-Q1"  (currently 3rd sheet)

# ===========================================================================
# Step 1 - "总计" sheet: add the 2022-Q4 summary row, shifting Q3/Q1 down.
# ===========================================================================

# Push the current row 3 (2022-Q1) down into row 4, then the current row 2
# (2022-Q3) down into row 3 - both via same-sheet copy so the existing cell
# formatting (style) travels with the data.
$summary.Range("A3:D3").Copy($summary.Range("A4:D4"))
$summary.Range("A2:D2").Copy($summary.Range("A3:D3"))

# Row 2 becomes the brand new 2022-Q4 entry.
$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 3
$summary.Range("D2").Value = 0.11

# Keep the running index column (A) sequential top to bottom.
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2

# ===========================================================================
# Step 2 - create the "2022-Q4" worksheet right after "总计".
# ===========================================================================

# Duplicate "2022-Q3" (same column layout/formatting) and drop the copy in
# right after "总计"; rename it to "2022-Q4" and overwrite its contents.
$q3.Copy($null, $summary)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# The template sheet only has one data row (row 2) - stamp it out twice more
# so we have three rows, each inheriting the same cell formatting.
$q4.Range("A2:H2").Copy($q4.Range("A3:H3"))
$q4.Range("A2:H2").Copy($q4.Range("A4:H4"))

# Row 2 - 011815 恒越优势精选混合
$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "'011815"
$q4.Range("C2").Value = "恒越优势精选混合"
$q4.Range("D2").Value = "'2.64"
$q4.Range("E2").Value = "'92.01"
$q4.Range("F2").Value = "'2.56"
$q4.Range("G2").Value = "'0.0676"
$q4.Range("H2").Value = 4

# Row 3 - 013028 恒越品质生活混合
$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "'013028"
$q4.Range("C3").Value = "恒越品质生活混合"
$q4.Range("D3").Value = "'1.25"
$q4.Range("E3").Value = "'90.58"
$q4.Range("F3").Value = "'2.36"
$q4.Range("G3").Value = "'0.0295"
$q4.Range("H3").Value = 9

# Row 4 - 159628 万家国证2000ETF
$q4.Range("A4").Value = 2
$q4.Range("B4").Value = "'159628"
$q4.Range("C4").Value = "万家国证2000ETF"
$q4.Range("D4").Value = "'2.55"
$q4.Range("E4").Value = "'97.28"
$q4.Range("F4").Value = "'0.48"
$q4.Range("G4").Value = "'0.0122"
$q4.Range("H4").Value = 3

# ===========================================================================
# Step 3 - restore the originally-selected tab (2022-Q1, the last sheet).
# ===========================================================================
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Activate()
